$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2035, 0.8099999999999999, 0.15, 0.04),
    @(2040, 0.8099999999999999, 0.15, 0.04),
    @(2045, 0.8099999999999999, 0.15, 0.04),
    @(2050, 0.8099999999999999, 0.15, 0.04)
)

$row = 5
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
